$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 22

    $ws.Cells.Item($row, 1).Value = 21                      # A: Trade #

    # Column B holds a date-looking string ("2026-02-16"). The source
    # workbook stores it as plain text, not an Excel date serial, so force
    # the cell to Text before writing it and then drop back to the Normal
    # style so no stray number format sticks around on the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"            # B: Date
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "22:59:28"              # C: Time
    $ws.Cells.Item($row, 4).Value = "base_strategy"         # D: Strategy
    $ws.Cells.Item($row, 5).Value = "UP"                    # E: Side
    $ws.Cells.Item($row, 6).Value = 0.5                     # F: Entry Price

    # Columns G and P are present but empty text cells in the source file
    # (not simply missing). A plain "" assignment clears/removes the cell
    # instead of leaving an empty text cell behind, so write a lone quote
    # prefix (forces Text, empty display value) and then drop the style
    # back to Normal so no stray quote-prefix formatting sticks around.
    $ws.Cells.Item($row, 7).Value = "'"                     # G: Exit Price
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "OPEN"                  # H: Status
    $ws.Cells.Item($row, 9).Value = 0                       # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0                      # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100                    # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                      # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                      # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                    # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason

    $ws.Cells.Item($row, 16).Value = "'"                    # P: Exit Reason
    $ws.Cells.Item($row, 16).Style = "Normal"

    $ws.Cells.Item($row, 17).Value = 0                      # Q: Duration (min)
}
